# Update countries & provincias Spain
# Refresh the "Pais" COVID-19 stats sheet with the 21:35 data snapshot
# (previous snapshot was 20:18): updated totals for several countries and
# two re-rankings (Guatemala overtakes Rumania; Islas Malvinas/Groenlandia
# swap position while tied).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 21 de Julio de 2020 a las 21:35"

# Helper: write a full data row (B..H) in one shot, leaving F (Casos criticos) as-is
function Set-CountryRow($row, $total, $nuevos, $activos, $recuperados, $muertesHoy, $muertes) {
    $ws.Range("B$row").Value = $total
    $ws.Range("C$row").Value = $nuevos
    $ws.Range("D$row").Value = $activos
    $ws.Range("E$row").Value = $recuperados
    $ws.Range("G$row").Value = $muertesHoy
    $ws.Range("H$row").Value = $muertes
}

# Row 4 - Estados Unidos
Set-CountryRow 4 3995841 34412 1864343 1986993 671 144505

# Row 6 - India
Set-CountryRow 6 1194085 39168 752393 412921 672 28771

# Row 21 - Alemania
Set-CountryRow 21 203846 359 188100 6566 7 9180

# Row 22 - Francia
Set-CountryRow 22 177338 584 79734 67439 13 30165

# Rows 48/49 - Guatemala overtakes Rumania in the ranking
$ws.Range("A48").Value = "Guatemala"
Set-CountryRow 48 40229 1190 26685 12013 29 1531

$ws.Range("A49").Value = "Rumania"
Set-CountryRow 49 39133 994 24454 12605 36 2074

# Row 77 - Costa Rica
Set-CountryRow 77 11811 277 3194 8549 2 68

# Row 92 - Guayana Francesa
Set-CountryRow 92 6851 106 4996 1816 1 39

# Rows 210/211 - Islas Malvinas / Groenlandia swap position (tied values)
$ws.Range("A210").Value = "Islas Malvinas"
$ws.Range("A211").Value = "Groenlandia"
